$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 183
$ws.Range("I11").Value = 183
$ws.Range("K11").Value = 183
$ws.Range("M11").Value = -43
$ws.Range("H64").Value = 24394258
$ws.Range("I64").Value = 3779.3428
$ws.Range("J64").Value = 166672050
$ws.Range("K64").Value = 3779.3428
$ws.Range("L64").Value = 166672050
$ws.Range("M64").Value = -3531.3428
$ws.Range("N64").Value = -166672546
$ws.Range("H67").Value = 24394258
$ws.Range("I67").Value = 3779.3428
$ws.Range("J67").Value = 166672050
$ws.Range("K67").Value = 3779.3428
$ws.Range("L67").Value = 166672050
$ws.Range("M67").Value = -2921.3428
$ws.Range("N67").Value = -166673766
$ws.Range("H86").Value = 1000000000
$ws.Range("I86").Value = 1000000000
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 1000000000
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -999998877
$ws.Range("N86").ClearContents()
$ws.Range("H88").Value = 11112294
$ws.Range("I88").Value = 25001012
$ws.Range("J88").Value = 1319.4
$ws.Range("K88").Value = 25001012
$ws.Range("L88").Value = 1319.4
$ws.Range("M88").Value = -25000606
$ws.Range("N88").Value = -2131.4
$ws.Range("H89").Value = 1000000000
$ws.Range("I89").Value = 1000000000
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 5000000000
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -4999994384
$ws.Range("N89").ClearContents()
$ws.Range("H91").Value = 11112294
$ws.Range("I91").Value = 25001012
$ws.Range("J91").Value = 1319.4
$ws.Range("K91").Value = 25001012
$ws.Range("L91").Value = 1319.4
$ws.Range("M91").Value = -24999608
$ws.Range("N91").Value = -4127.4
$ws.Range("H132").Value = 6389.5312
$ws.Range("I132").Value = 7391.25
$ws.Range("J132").Value = 3384.375
$ws.Range("K132").Value = 22173.75
$ws.Range("L132").Value = 10153.125
$ws.Range("M132").Value = -19643.75
$ws.Range("N132").Value = -15213.125
$ws.Range("H133").Value = 105000
$ws.Range("J133").Value = 105000
$ws.Range("L133").Value = 105000
$ws.Range("N133").Value = -115120
$ws.Range("H137").Value = 1115586
$ws.Range("I137").Value = 2381979.8
$ws.Range("K137").Value = 7145939.399999999
$ws.Range("M137").Value = -7143389.399999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8548.666999999999
$ws.Range("I32").Value = 8306.777
$ws.Range("J32").Value = 10000
$ws.Range("K32").Value = 8306.777
$ws.Range("L32").Value = 10000
$ws.Range("M32").Value = -8019.777
$ws.Range("N32").Value = -10574
$ws.Range("H61").Value = 2227.075
$ws.Range("I61").Value = 1729.3334
$ws.Range("J61").Value = 4573.5713
$ws.Range("K61").Value = 1729.3334
$ws.Range("L61").Value = 4573.5713
$ws.Range("M61").Value = -1517.3334
$ws.Range("N61").Value = -4997.5713
$ws.Range("H87").Value = 195000
$ws.Range("J87").Value = 195000
$ws.Range("L87").Value = 195000
$ws.Range("N87").Value = -197496
$ws.Range("H90").Value = 195000
$ws.Range("J90").Value = 195000
$ws.Range("L90").Value = 585000
$ws.Range("N90").Value = -597480
$ws.Range("H132").Value = 2428
$ws.Range("I132").Value = 2280.138
$ws.Range("J132").Value = 3500
$ws.Range("K132").Value = 6840.414
$ws.Range("L132").Value = 10500
$ws.Range("M132").Value = -4310.414
$ws.Range("N132").Value = -15560
$ws.Range("H136").Value = 2227.075
$ws.Range("I136").Value = 1729.3334
$ws.Range("J136").Value = 4573.5713
$ws.Range("K136").Value = 5188.0002
$ws.Range("L136").Value = 13720.7139
$ws.Range("M136").Value = -2638.0002
$ws.Range("N136").Value = -18820.7139

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1333
$ws.Range("I107").Value = 1109.3636
$ws.Range("K107").Value = 1109.3636
$ws.Range("M107").Value = 810.6364000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3949.743
$ws.Range("I31").Value = 2853.04
$ws.Range("J31").Value = 6691.5
$ws.Range("K31").Value = 2853.04
$ws.Range("L31").Value = 6691.5
$ws.Range("M31").Value = -2558.04
$ws.Range("N31").Value = -7281.5
$ws.Range("H34").Value = 3949.743
$ws.Range("I34").Value = 2853.04
$ws.Range("J34").Value = 6691.5
$ws.Range("K34").Value = 2853.04
$ws.Range("L34").Value = 6691.5
$ws.Range("M34").Value = -2651.04
$ws.Range("N34").Value = -7095.5
$ws.Range("H58").Value = 4273
$ws.Range("I58").Value = 3912
$ws.Range("J58").Value = 4333.1665
$ws.Range("K58").Value = 3912
$ws.Range("L58").Value = 4333.1665
$ws.Range("M58").Value = -3709
$ws.Range("N58").Value = -4739.1665
$ws.Range("H132").Value = 27782122
$ws.Range("I132").Value = 50002920
$ws.Range("J132").Value = 6121.75
$ws.Range("K132").Value = 150008760
$ws.Range("L132").Value = 18365.25
$ws.Range("M132").Value = -150006230
$ws.Range("N132").Value = -23425.25
$ws.Range("H134").Value = 3708.7646
$ws.Range("I134").Value = 3754.0833
$ws.Range("J134").Value = 3600
$ws.Range("K134").Value = 11262.2499
$ws.Range("L134").Value = 10800
$ws.Range("M134").Value = -8727.249899999999
$ws.Range("N134").Value = -15870
$ws.Range("H136").Value = 4273
$ws.Range("I136").Value = 3912
$ws.Range("J136").Value = 4333.1665
$ws.Range("K136").Value = 11736
$ws.Range("L136").Value = 12999.4995
$ws.Range("M136").Value = -9186
$ws.Range("N136").Value = -18099.4995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 157.63637
$ws.Range("I40").Value = 157.63637
$ws.Range("K40").Value = 630.54548
$ws.Range("M40").Value = -561.54548
$ws.Range("H107").Value = 517.8889
$ws.Range("J107").Value = 446.83334
$ws.Range("L107").Value = 1340.50002
$ws.Range("N107").Value = -5180.500019999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 111114820
$ws.Range("I80").Value = 200002510
$ws.Range("J80").Value = 5187
$ws.Range("K80").Value = 200002510
$ws.Range("L80").Value = 5187
$ws.Range("M80").Value = -200001512
$ws.Range("N80").Value = -7183
$ws.Range("H83").Value = 111114820
$ws.Range("I83").Value = 200002510
$ws.Range("J83").Value = 5187
$ws.Range("K83").Value = 1000012550
$ws.Range("L83").Value = 25935
$ws.Range("M83").Value = -1000007558
$ws.Range("N83").Value = -35919
$ws.Range("H132").Value = 2477.8096
$ws.Range("I132").Value = 2119.647
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 6358.941
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -3828.941
$ws.Range("N132").Value = -17060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1969.6
$ws.Range("I46").Value = 2122.5454
$ws.Range("K46").Value = 2122.5454
$ws.Range("M46").Value = -1934.5454
$ws.Range("H99").Value = 22421.143
$ws.Range("I99").Value = 22421.143
$ws.Range("K99").Value = 22421.143
$ws.Range("M99").Value = -19426.143
$ws.Range("I100").Value = 1299.2
$ws.Range("J100").Value = 1549
$ws.Range("K100").Value = 1299.2
$ws.Range("L100").Value = 1549
$ws.Range("M100").Value = -758.2
$ws.Range("N100").Value = -2631
$ws.Range("H132").Value = 4793.9614
$ws.Range("I132").Value = 2861.2354
$ws.Range("J132").Value = 8444.666999999999
$ws.Range("K132").Value = 8583.706200000001
$ws.Range("L132").Value = 25334.001
$ws.Range("M132").Value = -6053.706200000001
$ws.Range("N132").Value = -30394.001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H114").Value = 89499
$ws.Range("J114").Value = 89499
$ws.Range("L114").Value = 89499
$ws.Range("N114").Value = -98177
$ws.Range("H136").Value = 9021.432000000001
$ws.Range("J136").Value = 9939.308000000001
$ws.Range("L136").Value = 29817.924
$ws.Range("N136").Value = -34917.924
